$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Terra Mystica's "time" value (row 4, column D) from 5 to 30
$ws.Range("D4").Value = 30

# Update Tesla Vs Edison's "time" value (row 8, column D) from 6 to 20
$ws.Range("D8").Value = 20

# Add a new row for "Caverna"
$ws.Range("A9").Value = "Caverna"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 210
$ws.Range("E9").Value = $false

# Update selection to reflect the saved view state (D8 selected)
$ws.Range("D8").Select()

# Adjust the workbook window position/size as captured in the saved view
$excel.ActiveWindow.Left = -28920
$excel.ActiveWindow.Top = 1185
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15840
